# Adding selling and buying costs
# The RESOURCES sheet currently has a "costs_kWh" column (buying/operation cost).
# This change inserts a new "costs_sell_kWh" column right after it (selling cost),
# duplicating the buying-cost values, and documents it with a header comment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RESOURCES")

# Insert a new column F (pushing the old "reference" column from F to G).
$ws.Columns.Item(6).Insert()

# Header for the new column.
$ws.Cells.Item(1, 6).Value = "costs_sell_kWh"

# Selling-cost values mirror the existing buying-cost ("costs_kWh", column E) values.
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(3, 6).Formula = "=4.94/293"
$ws.Cells.Item(4, 6).Formula = "=0.2*0.75"
$ws.Cells.Item(5, 6).Value = 0.0001

# Document the new column, just like the other headers in this sheet are documented.
$ws.Cells.Item(1, 6).AddComment("Selling price in US`$(2015)/kWh(resource [thermal in case of fuels]).yr")

# Match the new active selection on the sheet.
$ws.Range("F1").Select() | Out-Null
